## add major enhancement to ppp
##
## 1) Refresh the cached "datetimeFigureOut" footer field (10/11/2018 -> 12/11/2018)
##    on the slide master and every slide layout.
## 2) Reposition/resize the "Straight Connector 18" connector (and flip it
##    horizontally) on slide 1.
## 3) Reposition/resize the "TextBox 1" (the small "X" cross-out mark) on slide 1.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached date-field text everywhere it is rendered from.
# ---------------------------------------------------------------------------
$newDate = "12/11/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every layout that hangs off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# 2) "Straight Connector 18" - flip horizontally + move/resize.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

$connector = $slide.Shapes.Item("Straight Connector 18")
$connector.Flip(0)  # msoFlipHorizontal
$connector.Left = 433.0250787401575
$connector.Top = 168.94287401574803
$connector.Width = 0.8591732283464567
$connector.Height = 219.92917322834646

# ---------------------------------------------------------------------------
# 3) "TextBox 1" (the small X cross-out mark) - move/resize.
# ---------------------------------------------------------------------------
$crossBox = $slide.Shapes.Item("TextBox 1")
$crossBox.Left = 419.3801968503937
$crossBox.Top = 372.39161417322833
$crossBox.Width = 6.040196850393701
$crossBox.Height = 46.04531496062992
